$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(48, "شادن", "ق", $false, 0),
    @(49, "شادن", "ا", $false, 0),
    @(50, "شادن", "ف", $false, 0),
    @(51, "شادن", "ز", $false, 0),
    @(52, "شادن", "ف", $true, 3.29),
    @(53, "شادن", "ف", $true, 0),
    @(54, "شادن", "د", $false, 0),
    @(55, "شادن", "ف", $true, 3),
    @(56, "شادن", "ا", $false, 0),
    @(57, "شادن", "د", $true, 0.44),
    @(58, "ٍshaden", "س", $false, 0),
    @(59, "shaden", "س", $false, 0),
    @(60, "shaden", "ه", $false, 0)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}
